$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of employee data (row 5)
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Mari"
$ws.Cells.Item(5, 3).Value = "Carolina"
$ws.Cells.Item(5, 4).Value = "Peroz"
$ws.Cells.Item(5, 5).Value = "Hans"
$ws.Cells.Item(5, 6).Value = "12344ee"
$ws.Cells.Item(5, 7).Value = "joakkkdjdd"
$ws.Cells.Item(5, 8).Value = "afafafafaf"

# fecha_nacimiento - keep as plain text, not an Excel date
$ws.Cells.Item(5, 9).NumberFormat = "@"
$ws.Cells.Item(5, 9).Value = "1985-07-11"
$ws.Cells.Item(5, 9).Style = "Normal"

$ws.Cells.Item(5, 10).Value = 39
$ws.Cells.Item(5, 11).Value = "Femenino"
$ws.Cells.Item(5, 12).Value = "Soltero"
$ws.Cells.Item(5, 13).Value = "DWd"
$ws.Cells.Item(5, 14).Value = "hdfhf"

# fecha_ingreso - keep as plain text, not an Excel date
$ws.Cells.Item(5, 15).NumberFormat = "@"
$ws.Cells.Item(5, 15).Value = "12-14-2000"
$ws.Cells.Item(5, 15).Style = "Normal"

$ws.Cells.Item(5, 16).Value = "Admon"
$ws.Cells.Item(5, 17).Value = "Quincenal"
$ws.Cells.Item(5, 18).Value = "Activo"
$ws.Cells.Item(5, 19).Value = "Banesco"

# numero_cuenta - keep as plain text so the long digit string is not
# coerced into a floating point number (which would lose precision)
$ws.Cells.Item(5, 20).NumberFormat = "@"
$ws.Cells.Item(5, 20).Value = "132242535353543535"
$ws.Cells.Item(5, 20).Style = "Normal"
